{"js": "// Highlight the key phrases in the two problem-statement bullets (yellow),\n// splitting each run into before/highlighted/after pieces while keeping the\n// existing italic red Times New Roman formatting intact.\n\nconst target1 = \"no proper method\";\nconst target2 = \"generating new plausible scenarios\";\n\nconst results1 = context.document.body.search(target1, { matchCase: true });\nresults1.load(\"items\");\nconst results2 = context.document.body.search(target2, { matchCase: true });\nresults2.load(\"items\");\nawait context.sync();\n\nif (results1.items.length > 0) {\n  results1.items[0].font.highlightColor = \"Yellow\";\n}\nif (results2.items.length > 0) {\n  results2.items[0].font.highlightColor = \"Yellow\";\n}\nawait context.sync();\n", "ps1": "# Highlight the key phrases in the two problem-statement bullets (yellow),\n# splitting each run into before/highlighted/after pieces while keeping the\n# existing italic red Times New Roman formatting intact.\n\n$d = $word.ActiveDocument\n\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Text = \"no proper method\"\n$rng1.Find.Execute() | Out-Null\n$rng1.Font.HighlightColorIndex = 7\n\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"generating new plausible scenarios\"\n$rng2.Find.Execute() | Out-Null\n$rng2.Font.HighlightColorIndex = 7\n"}
